# Automatische test-sync: 2025-06-18 09:00:10
# Adds a new "Afmelding nieuwsbrief" (newsletter unsubscribe) mail-log entry
# on the Logs sheet, and bumps the matching category count on the
# Dashboard sheet (which drives the bar chart).

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: append row 3 -------------------------------------------------
$logs.Range("A3").Value = "Afmelding nieuwsbrief"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("C3").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Range("D3").Value = "Afmelding"
$logs.Range("F3").Value = "2025-06-18 08:30:11"
$logs.Range("G3").Value = "Nee"

# Extend the conditional formatting ranges so the new row is covered too.
$dCond = $logs.Range("D2").FormatConditions.Item(1)
$dCond.ModifyAppliesToRange($logs.Range("D2:D3"))

$gCond = $logs.Range("G2").FormatConditions.Item(1)
$gCond.ModifyAppliesToRange($logs.Range("G2:G3"))

# --- Dashboard sheet: append row 3 --------------------------------------------
$dashboard.Range("A3").Value = "Afmelding"
$dashboard.Range("B3").Value = 1

# --- Chart: extend the category/value series references to include row 3 -----
$chartObj = $dashboard.ChartObjects(1)
$series = $chartObj.Chart.SeriesCollection(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$3,Dashboard!`$B`$2:`$B`$3,1)"
